$wb = $excel.ActiveWorkbook

# --- Sheet "k" (sheet1) ---
$wsK = $wb.Worksheets.Item("k")

# Move the NC_000913 / NC_010468 labels from A16:B16 to I4:J4
$wsK.Range("A16").Value = $null
$wsK.Range("B16").Value = $null
$wsK.Range("I4").Value = "NC_000913"
$wsK.Range("J4").Value = "NC_010468"

# Update the selection shown on this sheet
$wsK.Range("I4:J4").Select()

# --- Sheet "q" (sheet4) ---
$wsQ = $wb.Worksheets.Item("q")

# Update the selection shown on this sheet
$wsQ.Range("B2:B11").Select()

# Keep sheet "k" as the active/displayed tab (as in the original workbook)
$wsK.Activate()
